$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''67.429.86'
$ws.Range("E2").Value = '  -2.86%  '

$ws.Range("D3").Value = '''3.516.84'
$ws.Range("E3").Value = '  -3.97%  '

$ws.Range("E4").Value = '  +0.16%  '

$ws.Range("D5").Value = '''616.72'
$ws.Range("E5").Value = '  -4.39%  '

$ws.Range("D6").Value = '''153.39'
$ws.Range("E6").Value = '  -3.61%  '

$ws.Range("D7").Value = '''3.512.36'
$ws.Range("E7").Value = '  -3.99%  '

$ws.Range("E8").Value = '  +0.20%  '

$ws.Range("E9").Value = '  -2.10%  '

$ws.Range("E10").Value = '  -2.39%  '

$ws.Range("D11").Value = '''6.89'
$ws.Range("E11").Value = '  -2.39%  '

$ws.Range("D12").Value = '''0.432'
$ws.Range("E12").Value = '  -1.36%  '

$ws.Range("E13").Value = '  -3.03%  '

$ws.Range("D14").Value = '''32.28'
$ws.Range("E14").Value = '  +0.18%  '

$ws.Range("D15").Value = '''4.096.35'
$ws.Range("E15").Value = '  -4.22%  '

$ws.Range("D16").Value = '''3.539.48'
$ws.Range("E16").Value = '  -3.20%  '

$ws.Range("D17").Value = '''67.524.10'
$ws.Range("E17").Value = '  -2.72%  '

$ws.Range("E18").Value = '  +0.66%  '

$ws.Range("B19").Value = 'Chainlink'
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D19").Value = '''15.57'
$ws.Range("E19").Value = '  -2.18%  '

$ws.Range("B20").Value = 'Polkadot'
$ws.Range("C20").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D20").Value = '''6.36'
$ws.Range("E20").Value = '  -0.70%  '

$ws.Range("D21").Value = '''453.73'
$ws.Range("E21").Value = '  -2.43%  '

$ws.Range("D22").Value = '''9.40'
$ws.Range("E22").Value = '  -3.16%  '

$ws.Range("D23").Value = '''0.644'
$ws.Range("E23").Value = '  +0.36%  '

$ws.Range("D24").Value = '''78.05'
$ws.Range("E24").Value = '  -1.72%  '

$ws.Range("D25").Value = '''3.678.58'
$ws.Range("E25").Value = '  -3.37%  '

$ws.Range("E26").Value = '  +0.09%  '

$ws.Range("E27").Value = '  -3.66%  '

$ws.Range("D28").Value = '''10.51'
$ws.Range("E28").Value = '  -1.61%  '

$ws.Range("D29").Value = '''8.36'
$ws.Range("E29").Value = '  -5.49%  '

$ws.Range("D30").Value = '''2.58'
$ws.Range("E30").Value = '  +0.06%  '

$ws.Range("D31").Value = '''1.69'
$ws.Range("E31").Value = '  +2.52%  '

$ws.Range("E32").Value = '  +0.03%  '

$ws.Range("D33").Value = '''25.98'
$ws.Range("E33").Value = '  -1.97%  '

$ws.Range("D34").Value = '''1.91'
$ws.Range("E34").Value = '  -3.24%  '

$ws.Range("D35").Value = '''6.23'
$ws.Range("E35").Value = '  -2.64%  '

$ws.Range("D36").Value = '''0.158'
$ws.Range("E36").Value = '  -2.77%  '

$ws.Range("D37").Value = '''3.535.68'
$ws.Range("E37").Value = '  -3.10%  '

$ws.Range("D38").Value = '''8.01'
$ws.Range("E38").Value = '  -3.60%  '

$ws.Range("E39").Value = '  +0.03%  '

$ws.Range("E40").Value = '  +0.12%  '

$ws.Range("D41").Value = '''176.46'
$ws.Range("E41").Value = '  -0.91%  '

$ws.Range("D42").Value = '''5.61'
$ws.Range("E42").Value = '  -4.56%  '

$ws.Range("D43").Value = '''0.0879'
$ws.Range("E43").Value = '  -1.02%  '

$ws.Range("D44").Value = '''2.10'
$ws.Range("E44").Value = '  -3.36%  '

$ws.Range("D45").Value = '''0.886'
$ws.Range("E45").Value = '  -4.29%  '

$ws.Range("D46").Value = '''29.27'
$ws.Range("E46").Value = '  +9.23%  '

$ws.Range("D47").Value = '''45.74'
$ws.Range("E47").Value = '  -1.82%  '

$ws.Range("D48").Value = '''2.59'
$ws.Range("E48").Value = '  -3.96%  '

$ws.Range("D49").Value = '''7.65'
$ws.Range("E49").Value = '  -1.52%  '

$ws.Range("E50").Value = '  -3.03%  '

$ws.Range("E51").Value = '  -2.80%  '
